# Update data: 2025-11-05 15:01
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump "Last Updated" timestamp by a minute ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 03:01 PM"

# --- Industry Analysis sheet: refresh the "1 Year" column (F2:F76) ---
$ws2 = $wb.Worksheets.Item("Industry Analysis")
$ws2.Cells.Item(2,6).Value = 21.0016
$ws2.Cells.Item(3,6).Value = -16.2396
$ws2.Cells.Item(4,6).Value = 27.1317
$ws2.Cells.Item(5,6).Value = -50.6494
$ws2.Cells.Item(6,6).Value = 53.2813
$ws2.Cells.Item(7,6).Value = -8.106199999999999
$ws2.Cells.Item(8,6).Value = -9.552099999999999
$ws2.Cells.Item(9,6).Value = 36.3756
$ws2.Cells.Item(10,6).Value = -6.1314
$ws2.Cells.Item(11,6).Value = 31.9081
$ws2.Cells.Item(12,6).Value = -18.4955
$ws2.Cells.Item(13,6).Value = 14.0155
$ws2.Cells.Item(14,6).Value = -36.0718
$ws2.Cells.Item(15,6).Value = -0.1622
$ws2.Cells.Item(16,6).Value = 0.1459
$ws2.Cells.Item(17,6).Value = -22.0012
$ws2.Cells.Item(18,6).Value = 1.0561
$ws2.Cells.Item(19,6).Value = -27.708
$ws2.Cells.Item(20,6).Value = 47.7309
$ws2.Cells.Item(21,6).Value = 12.0959
$ws2.Cells.Item(22,6).Value = 95.1491
$ws2.Cells.Item(23,6).Value = -50.2657
$ws2.Cells.Item(24,6).Value = -13.3427
$ws2.Cells.Item(25,6).Value = -9.9316
$ws2.Cells.Item(26,6).Value = 5.8244
$ws2.Cells.Item(27,6).Value = -32.7692
$ws2.Cells.Item(28,6).Value = -24.8224
$ws2.Cells.Item(29,6).Value = -18.4191
$ws2.Cells.Item(30,6).Value = 25.8569
$ws2.Cells.Item(31,6).Value = 58.4712
$ws2.Cells.Item(32,6).Value = -3.3862
$ws2.Cells.Item(33,6).Value = -6.3282
$ws2.Cells.Item(34,6).Value = 27.7203
$ws2.Cells.Item(35,6).Value = 4.4873
$ws2.Cells.Item(36,6).Value = -4.9458
$ws2.Cells.Item(37,6).Value = 3.6074
$ws2.Cells.Item(38,6).Value = -23.3973
$ws2.Cells.Item(39,6).Value = 8.7355
$ws2.Cells.Item(40,6).Value = -5.8541
$ws2.Cells.Item(41,6).Value = -8.3934
$ws2.Cells.Item(42,6).Value = 20.3818
$ws2.Cells.Item(43,6).Value = 14.3164
$ws2.Cells.Item(44,6).Value = -12.6846
$ws2.Cells.Item(45,6).Value = 28.4075
$ws2.Cells.Item(46,6).Value = -1.1135
$ws2.Cells.Item(47,6).Value = -37.1997
$ws2.Cells.Item(48,6).Value = -29.8569
$ws2.Cells.Item(49,6).Value = -27.5511
$ws2.Cells.Item(50,6).Value = -49.7478
$ws2.Cells.Item(51,6).Value = -51.8002
$ws2.Cells.Item(52,6).Value = -38.5254
$ws2.Cells.Item(53,6).Value = -12.4886
$ws2.Cells.Item(54,6).Value = -5.0725
$ws2.Cells.Item(55,6).Value = -17.7445
$ws2.Cells.Item(56,6).Value = -26.636
$ws2.Cells.Item(57,6).Value = -29.3361
$ws2.Cells.Item(58,6).Value = -11.9574
$ws2.Cells.Item(59,6).Value = -24.5687
$ws2.Cells.Item(60,6).Value = -12.3
$ws2.Cells.Item(61,6).Value = -10.9446
$ws2.Cells.Item(62,6).Value = -17.1229
$ws2.Cells.Item(63,6).Value = -9.5038
$ws2.Cells.Item(64,6).Value = 54.2749
$ws2.Cells.Item(65,6).Value = -43.4736
$ws2.Cells.Item(66,6).Value = 13.2687
$ws2.Cells.Item(67,6).Value = 12.7149
$ws2.Cells.Item(68,6).Value = 24.8057
$ws2.Cells.Item(69,6).Value = -17.0328
$ws2.Cells.Item(70,6).Value = -6.8927
$ws2.Cells.Item(71,6).Value = 13.6034
$ws2.Cells.Item(72,6).Value = 3.9995
$ws2.Cells.Item(73,6).Value = -16.226
$ws2.Cells.Item(74,6).Value = -16.2448
$ws2.Cells.Item(75,6).Value = 28.6924
$ws2.Cells.Item(76,6).Value = 48.9752

# --- Stock List sheet: new top row (CAPTRU-RE1) pushed in, oldest row dropped ---
$ws3 = $wb.Worksheets.Item("Stock List")
$ws3.Rows.Item(2).Insert()
$ws3.Rows.Item(77).Delete()
$ws3.Rows.Item(2).ClearFormats()
$ws3.Cells.Item(2,1).Value = [char]0x1F4CB
$ws3.Cells.Item(2,2).Value = "CAPTRU-RE1"
$ws3.Cells.Item(2,3).Value = "CAPTRU-RE1"
$ws3.Cells.Item(2,4).Value = 5.67
$ws3.Cells.Item(2,5).Value = -11.9565
$ws3.Cells.Item(2,6).Value = "N/A"
$ws3.Cells.Item(2,7).Value = "N/A"
$ws3.Cells.Item(2,8).Value = 0
